$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# --- Extend formatting of the new columns (H:N) to match existing header/body styles ---
# Header row 1: copy the format of the last existing header cell (G1) onto H1:N1
$ws.Range("G1").Copy()
$ws.Range("H1:N1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Data rows 2-6: copy the format of the last existing body cell in each row (G<row>) onto H<row>:N<row>
$ws.Range("G2:G6").Copy()
$ws.Range("H2:N6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

# --- Row 1 (header) ---
$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "debtor"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Row 2 ---
$ws.Range("B2").Value = "房屋貸款"
$ws.Range("C2").Value = "管碧玲"
$ws.Range("D2").Value = "合作金庫西門支庫臺北市萬華區昆明街"
$ws.Range("E2").Value = 908219
$ws.Range("F2").Value = "88年09‘月17日"
$ws.Range("G2").Value = "購買房屋"
$ws.Range("H2").Value = "debt"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2011-12-16"
$ws.Range("K2").Value = "管碧玲"
$ws.Range("L2").Value = 1374
$ws.Range("M2").Value = "tmp32301"
$ws.Range("N2").Value = 111

# --- Row 3 ---
$ws.Range("B3").Value = "房屋貸款"
$ws.Range("C3").Value = "許陽明"
$ws.Range("D3").Value = "高雄銀行三民分行高雄市三民區博愛一路"
$ws.Range("E3").Value = 7473463
$ws.Range("F3").Value = "99年09月17日"
$ws.Range("G3").Value = "購買房屋"
$ws.Range("H3").Value = "debt"
$ws.Range("I3").Value = "normal"
$ws.Range("J3").Value = "2011-12-16"
$ws.Range("K3").Value = "管碧玲"
$ws.Range("L3").Value = 1374
$ws.Range("M3").Value = "tmp32301"
$ws.Range("N3").Value = 112

# --- Row 4 ---
$ws.Range("B4").Value = "借款"
$ws.Range("C4").Value = "管碧玲"
$ws.Range("D4").Value = "管東隆彰化縣員林鎮林森路"
$ws.Range("E4").Value = 3000000
$ws.Range("F4").Value = "97年09月01曰"
$ws.Range("G4").Value = "借款"
$ws.Range("H4").Value = "debt"
$ws.Range("I4").Value = "normal"
$ws.Range("J4").Value = "2011-12-16"
$ws.Range("K4").Value = "管碧玲"
$ws.Range("L4").Value = 1374
$ws.Range("M4").Value = "tmp32301"
$ws.Range("N4").Value = 113

# --- Row 5 ---
$ws.Range("B5").Value = "借款"
$ws.Range("C5").Value = "管碧玲"
$ws.Range("D5").Value = "吳麗珠臺北市大安區和平東路"
$ws.Range("E5").Value = 700000
$ws.Range("F5").Value = "97年09月16日"
$ws.Range("G5").Value = "借款"
$ws.Range("H5").Value = "debt"
$ws.Range("I5").Value = "normal"
$ws.Range("J5").Value = "2011-12-16"
$ws.Range("K5").Value = "管碧玲"
$ws.Range("L5").Value = 1374
$ws.Range("M5").Value = "tmp32301"
$ws.Range("N5").Value = 114

# --- Row 6 ---
$ws.Range("B6").Value = "借款"
$ws.Range("C6").Value = "管碧玲"
$ws.Range("D6").Value = "管東隆彰化縣員林鎮林森路"
$ws.Range("E6").Value = 400000
$ws.Range("F6").Value = "98年07月07日"
$ws.Range("G6").Value = "借款"
$ws.Range("H6").Value = "debt"
$ws.Range("I6").Value = "normal"
$ws.Range("J6").Value = "2011-12-16"
$ws.Range("K6").Value = "管碧玲"
$ws.Range("L6").Value = 1374
$ws.Range("M6").Value = "tmp32301"
$ws.Range("N6").Value = 115
